$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.443.97"
$ws.Range("D3").Value = "2.511.90"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "2.511.63"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  +4.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "2.975.65"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "69.303.05"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "2.523.70"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").Value = "2.658.91"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "0.0₃0883"
$ws.Range("E30").Value = "  -3.47%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "460.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.79%  "
$ws.Range("E33").Value = "  -4.62%  "
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  -7.36%  "
$ws.Range("E46").Value = "  -7.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.518"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -3.37%  "
